$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.268868803977966
$ws.Range("B1").Value = 1.435236692428589
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.036584615707397
$ws.Range("E1").Value = 0.8921931385993958
